$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "53.658.41"
Set-TextValue $ws "E2" "  -5.01%  "
Set-TextValue $ws "D3" "2.211.91"
Set-TextValue $ws "E3" "  -6.95%  "
Set-TextValue $ws "D4" "1.01"
Set-TextValue $ws "E4" "  +0.46%  "
Set-TextValue $ws "D5" "487.30"
Set-TextValue $ws "E5" "  -3.78%  "
Set-TextValue $ws "D6" "125.46"
Set-TextValue $ws "E6" "  -4.01%  "
Set-TextValue $ws "D7" "1.00"
Set-TextValue $ws "E7" "  +0.15%  "
Set-TextValue $ws "E8" "  -4.13%  "
Set-TextValue $ws "D9" "2.240.59"
Set-TextValue $ws "E9" "  -6.26%  "
Set-TextValue $ws "E10" "  -6.52%  "
Set-TextValue $ws "E11" "  -0.16%  "
Set-TextValue $ws "D12" "0.321"
Set-TextValue $ws "E12" "  -2.95%  "
Set-TextValue $ws "E13" "  -4.76%  "
Set-TextValue $ws "D14" "2.609.27"
Set-TextValue $ws "E14" "  -6.78%  "
Set-TextValue $ws "D15" "21.25"
Set-TextValue $ws "E15" "  -1.79%  "
Set-TextValue $ws "D16" "53.611.97"
Set-TextValue $ws "E16" "  -5.03%  "
Set-TextValue $ws "E17" "  -4.00%  "
Set-TextValue $ws "D18" "2.257.12"
Set-TextValue $ws "E18" "  -6.14%  "
Set-TextValue $ws "D19" "9.67"
Set-TextValue $ws "E19" "  -3.93%  "
Set-TextValue $ws "E20" "  -1.91%  "
Set-TextValue $ws "D21" "296.46"
Set-TextValue $ws "E21" "  -4.24%  "
Set-TextValue $ws "D22" "6.18"
Set-TextValue $ws "E22" "  -1.67%  "
Set-TextValue $ws "E23" "  -0.26%  "
Set-TextValue $ws "D24" "63.90"
Set-TextValue $ws "E24" "  -3.55%  "
Set-TextValue $ws "D25" "1.00"
Set-TextValue $ws "E25" "  +0.33%  "
Set-TextValue $ws "D26" "0.368"
Set-TextValue $ws "E26" "  -1.09%  "
Set-TextValue $ws "E27" "  -1.17%  "
Set-TextValue $ws "D28" "2.314.32"
Set-TextValue $ws "E28" "  -7.07%  "
Set-TextValue $ws "D29" "7.04"
Set-TextValue $ws "E29" "  -3.14%  "
Set-TextValue $ws "D30" "162.78"
Set-TextValue $ws "E30" "  -6.30%  "
Set-TextValue $ws "E31" "  -4.03%  "
Set-TextValue $ws "D32" "0.999"
Set-TextValue $ws "E32" "  +0.01%  "
Set-TextValue $ws "E33" "  -0.96%  "
Set-TextValue $ws "D34" "0.0₃0669"
Set-TextValue $ws "E34" "  -6.31%  "
Set-TextValue $ws "D35" "0.992"
Set-TextValue $ws "E35" "  -0.35%  "
Set-TextValue $ws "D36" "1.06"
Set-TextValue $ws "E36" "  -1.72%  "
Set-TextValue $ws "D37" "17.29"
Set-TextValue $ws "E37" "  -2.16%  "
Set-TextValue $ws "E38" "  -1.08%  "
Set-TextValue $ws "D39" "0.832"
Set-TextValue $ws "E39" "  +0.62%  "
Set-TextValue $ws "E40" "  -3.74%  "
Set-TextValue $ws "D41" "35.19"
Set-TextValue $ws "E41" "  -3.62%  "
Set-TextValue $ws "D43" "1.38"
Set-TextValue $ws "E43" "  -1.19%  "
Set-TextValue $ws "D44" "127.92"
Set-TextValue $ws "E44" "  +0.37%  "
Set-TextValue $ws "E45" "  -2.45%  "
Set-TextValue $ws "D46" "4.78"
Set-TextValue $ws "E46" "  +0.34%  "
Set-TextValue $ws "D47" "0.0884"
Set-TextValue $ws "E47" "  -1.63%  "
Set-TextValue $ws "E48" "  -5.05%  "
Set-TextValue $ws "D49" "235.22"
Set-TextValue $ws "E49" "  -2.11%  "
Set-TextValue $ws "E50" "  -2.21%  "
Set-TextValue $ws "E51" "  -3.27%  "
